$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 150, shifting existing rows 150-152 down to 151-153.
$ws.Rows.Item(150).Insert()

# Populate the newly inserted row 150 with the new record.
$ws.Cells.Item(150, 1).Value = 5
$ws.Cells.Item(150, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(150, 3).Value = "Maule"
$ws.Cells.Item(150, 4).Value = 44628
$ws.Cells.Item(150, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(150, 5).Value = 7
$ws.Cells.Item(150, 6).Value = 100112031
$ws.Cells.Item(150, 7).Value = "Poroto verde"
$ws.Cells.Item(150, 8).Value = "Sin especificar"
$ws.Cells.Item(150, 9).Value = "Primera"
$ws.Cells.Item(150, 10).Value = 200
$ws.Cells.Item(150, 11).Value = 30000
$ws.Cells.Item(150, 12).Value = 30000
$ws.Cells.Item(150, 13).Value = 30000
$ws.Cells.Item(150, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(150, 15).Value = "Región del Maule"
$ws.Cells.Item(150, 16).Value = 1200
$ws.Cells.Item(150, 17).Value = 25
$ws.Cells.Item(150, 18).Value = "Hortaliza"
